$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 2374.75
$ws.Range("J12").Value = 4499.25
$ws.Range("L12").Value = 4499.25
$ws.Range("N12").Value = -4839.25

$ws.Range("H64").Value = 9750.75
$ws.Range("J64").Value = 9667.666999999999
$ws.Range("L64").Value = 9667.666999999999
$ws.Range("N64").Value = -10163.667

$ws.Range("H67").Value = 9750.75
$ws.Range("J67").Value = 9667.666999999999
$ws.Range("L67").Value = 9667.666999999999
$ws.Range("N67").Value = -11383.667

$ws.Range("H70").Value = 7336
$ws.Range("I70").Value = 921.75
$ws.Range("J70").Value = 13750.25
$ws.Range("K70").Value = 2765.25
$ws.Range("L70").Value = 41250.75
$ws.Range("M70").Value = -2495.25
$ws.Range("N70").Value = -41790.75

$ws.Range("H73").Value = 7336
$ws.Range("I73").Value = 921.75
$ws.Range("J73").Value = 13750.25
$ws.Range("K73").Value = 2765.25
$ws.Range("L73").Value = 41250.75
$ws.Range("M73").Value = -1829.25
$ws.Range("N73").Value = -43122.75

$ws.Range("H132").Value = 1130.2325
$ws.Range("I132").Value = 1037.9524
$ws.Range("K132").Value = 3113.857199999999
$ws.Range("M132").Value = -583.8571999999995

$ws.Range("H138").Value = 3553.4546
$ws.Range("J138").Value = 4498.3335
$ws.Range("L138").Value = 13495.0005
$ws.Range("N138").Value = -23775.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5268.278
$ws.Range("I32").Value = 4080.3333
$ws.Range("K32").Value = 4080.3333
$ws.Range("M32").Value = -3793.3333

$ws.Range("H45").Value = 90911816
$ws.Range("I45").Value = 142858600
$ws.Range("J45").Value = 4970.75
$ws.Range("K45").Value = 142858600
$ws.Range("L45").Value = 4970.75
$ws.Range("M45").Value = -142858223
$ws.Range("N45").Value = -5724.75

$ws.Range("H61").Value = 5380.7437
$ws.Range("I61").Value = 4601.3516
$ws.Range("K61").Value = 4601.3516
$ws.Range("M61").Value = -4389.3516

$ws.Range("H97").Value = 2071.5715
$ws.Range("I97").Value = 3255.25
$ws.Range("J97").Value = 493.33334
$ws.Range("K97").Value = 3255.25
$ws.Range("L97").Value = 493.33334
$ws.Range("M97").Value = -2759.25
$ws.Range("N97").Value = -1485.33334

$ws.Range("H110").Value = 2580.923
$ws.Range("I110").Value = 1316.5555
$ws.Range("K110").Value = 1316.5555
$ws.Range("M110").Value = 728.4445000000001

$ws.Range("H136").Value = 5380.7437
$ws.Range("I136").Value = 4601.3516
$ws.Range("K136").Value = 13804.0548
$ws.Range("M136").Value = -11254.0548

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3497.7856
$ws.Range("I58").Value = 1788.0834
$ws.Range("K58").Value = 1788.0834
$ws.Range("M58").Value = -1585.0834

$ws.Range("H132").Value = 2755.6924
$ws.Range("I132").Value = 1364.375
$ws.Range("J132").Value = 4981.8
$ws.Range("K132").Value = 4093.125
$ws.Range("L132").Value = 14945.4
$ws.Range("M132").Value = -1563.125
$ws.Range("N132").Value = -20005.4

$ws.Range("H136").Value = 3497.7856
$ws.Range("I136").Value = 1788.0834
$ws.Range("K136").Value = 5364.2502
$ws.Range("M136").Value = -2814.2502

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 3.5
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 3.5
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 10.5
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -358.5

$ws.Range("H21").Value = 100
$ws.Range("J21").Value = 100
$ws.Range("L21").Value = 300
$ws.Range("N21").Value = -646

$ws.Range("H32").Value = 12616.667
$ws.Range("I32").Value = 250
$ws.Range("J32").Value = 22510
$ws.Range("K32").Value = 750
$ws.Range("L32").Value = 67530
$ws.Range("M32").Value = -467
$ws.Range("N32").Value = -68096

$ws.Range("H34").Value = 5985128
$ws.Range("J34").Value = 3703.75
$ws.Range("L34").Value = 11111.25
$ws.Range("N34").Value = -11279.25

$ws.Range("H38").Value = 48.625
$ws.Range("I38").Value = 95
$ws.Range("J38").Value = 20.8
$ws.Range("K38").Value = 285
$ws.Range("L38").Value = 62.40000000000001
$ws.Range("M38").Value = 62
$ws.Range("N38").Value = -756.4

$ws.Range("H39").Value = 2760
$ws.Range("J39").Value = 3550
$ws.Range("L39").Value = 10650
$ws.Range("N39").Value = -11238

$ws.Range("H60").Value = 37037972
$ws.Range("J60").Value = 2650
$ws.Range("L60").Value = 7950
$ws.Range("N60").Value = -8452

$ws.Range("H76").Value = 14009.667
$ws.Range("J76").Value = 19015
$ws.Range("L76").Value = 57045
$ws.Range("N76").Value = -57811

$ws.Range("H79").Value = 14009.667
$ws.Range("J79").Value = 19015
$ws.Range("L79").Value = 57045
$ws.Range("N79").Value = -59697

$ws.Range("H122").Value = 1619.5
$ws.Range("J122").Value = 1619.5
$ws.Range("L122").Value = 14575.5
$ws.Range("N122").Value = -19475.5

$ws.Range("H131").Value = 35495696
$ws.Range("I131").Value = 41670670
$ws.Range("J131").Value = 32408208
$ws.Range("K131").Value = 125012010
$ws.Range("L131").Value = 97224624
$ws.Range("M131").Value = -125006970
$ws.Range("N131").Value = -97234704

$ws.Range("H139").Value = 5068
$ws.Range("I139").Value = 1460.5
$ws.Range("J139").Value = 12283
$ws.Range("K139").Value = 4381.5
$ws.Range("L139").Value = 36849
$ws.Range("M139").Value = 758.5
$ws.Range("N139").Value = -47129

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 16569.572
$ws.Range("I70").Value = 11000
$ws.Range("J70").Value = 17497.834
$ws.Range("K70").Value = 11000
$ws.Range("L70").Value = 17497.834
$ws.Range("M70").Value = -10730
$ws.Range("N70").Value = -18037.834

$ws.Range("H73").Value = 16569.572
$ws.Range("I73").Value = 11000
$ws.Range("J73").Value = 17497.834
$ws.Range("K73").Value = 11000
$ws.Range("L73").Value = 17497.834
$ws.Range("M73").Value = -10064
$ws.Range("N73").Value = -19369.834

$ws.Range("H80").Value = 6197.25
$ws.Range("I80").Value = 4357.875
$ws.Range("K80").Value = 4357.875
$ws.Range("M80").Value = -3359.875

$ws.Range("H83").Value = 6197.25
$ws.Range("I83").Value = 4357.875
$ws.Range("K83").Value = 21789.375
$ws.Range("M83").Value = -16797.375

$ws.Range("H97").Value = 2039.7
$ws.Range("I97").Value = 1616.6
$ws.Range("K97").Value = 1616.6
$ws.Range("M97").Value = -1120.6

$ws.Range("H122").Value = 5398.1665
$ws.Range("J122").Value = 5761.6
$ws.Range("L122").Value = 17284.8
$ws.Range("N122").Value = -22184.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5202
$ws.Range("I22").Value = 1513.4286
$ws.Range("J22").Value = 7549.273
$ws.Range("K22").Value = 1513.4286
$ws.Range("L22").Value = 7549.273
$ws.Range("M22").Value = -1218.4286
$ws.Range("N22").Value = -8139.273

$ws.Range("H27").Value = 5202
$ws.Range("I27").Value = 1513.4286
$ws.Range("J27").Value = 7549.273
$ws.Range("K27").Value = 1513.4286
$ws.Range("L27").Value = 7549.273
$ws.Range("M27").Value = -1406.4286
$ws.Range("N27").Value = -7763.273

$ws.Range("H46").Value = 2069.2
$ws.Range("J46").Value = 2130.7368
$ws.Range("L46").Value = 2130.7368
$ws.Range("N46").Value = -2506.7368

$ws.Range("H82").Value = 3071.0386
$ws.Range("I82").Value = 2726.158
$ws.Range("J82").Value = 4007.1428
$ws.Range("K82").Value = 2726.158
$ws.Range("L82").Value = 4007.1428
$ws.Range("M82").Value = -2365.158
$ws.Range("N82").Value = -4729.1428

$ws.Range("H85").Value = 3071.0386
$ws.Range("I85").Value = 2726.158
$ws.Range("J85").Value = 4007.1428
$ws.Range("K85").Value = 2726.158
$ws.Range("L85").Value = 4007.1428
$ws.Range("M85").Value = -1478.158
$ws.Range("N85").Value = -6503.1428

$ws.Range("H122").Value = 3083.1667
$ws.Range("I122").Value = 3625
$ws.Range("J122").Value = 1999.5
$ws.Range("K122").Value = 10875
$ws.Range("L122").Value = 5998.5
$ws.Range("M122").Value = -8425
$ws.Range("N122").Value = -10898.5

$ws.Range("H132").Value = 2804.9062
$ws.Range("I132").Value = 1406.9584
$ws.Range("K132").Value = 4220.8752
$ws.Range("M132").Value = -1690.8752

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4853.957
$ws.Range("I132").Value = 4043.1562
$ws.Range("K132").Value = 12129.4686
$ws.Range("M132").Value = -9599.4686

$ws.Range("H136").Value = 5345.857
$ws.Range("I136").Value = 2698.2
$ws.Range("K136").Value = 8094.599999999999
$ws.Range("M136").Value = -5544.599999999999
